$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 47.86240033333333
$ws.Range("H2").Value = 143.587201
$ws.Range("I2").Value = 0.1228118231805696
$ws.Range("J2").Value = 0.1228118231805696
$ws.Range("M2").Value = 266.9240163333333
$ws.Range("N2").Value = 800.7720489999999
$ws.Range("O2").Value = 0.7873936103073201
$ws.Range("P2").Value = 0.78739361030732
$ws.Range("Q2").Value = 12775.6241283272
$ws.Range("R2").Value = 114980.6171549448
$ws.Range("S2").Value = 0.09670124484257295
$ws.Range("T2").Value = 0.09670124484257292
$ws.Range("G3").Value = 47.86240033333333
$ws.Range("H3").Value = 143.587201
$ws.Range("I3").Value = 0.1228118231805696
$ws.Range("J3").Value = 0.1228118231805696
$ws.Range("O3").Value = 0.1525285446808506
$ws.Range("P3").Value = 0.1525285446808506
$ws.Range("Q3").Value = 2474.807174168903
$ws.Range("R3").Value = 22273.26456752012
$ws.Range("S3").Value = 0.01873230865933424
$ws.Range("T3").Value = 0.01873230865933424
$ws.Range("G4").Value = 47.86240033333333
$ws.Range("H4").Value = 143.587201
$ws.Range("I4").Value = 0.1228118231805696
$ws.Range("J4").Value = 0.1228118231805696
$ws.Range("M4").Value = 7.729888333333332
$ws.Range("N4").Value = 23.189665
$ws.Range("O4").Value = 0.02280223700236483
$ws.Range("P4").Value = 0.02280223700236483
$ws.Range("Q4").Value = 369.9710099419627
$ws.Range("R4").Value = 3329.739089477664
$ws.Range("S4").Value = 0.002800384298855871
$ws.Range("T4").Value = 0.002800384298855871
$ws.Range("G5").Value = 47.86240033333333
$ws.Range("H5").Value = 143.587201
$ws.Range("I5").Value = 0.1228118231805696
$ws.Range("J5").Value = 0.1228118231805696
$ws.Range("M5").Value = 3.652135
$ws.Range("N5").Value = 10.956405
$ws.Range("O5").Value = 0.01077335716164485
$ws.Range("P5").Value = 0.01077335716164485
$ws.Range("Q5").Value = 174.7999474413783
$ws.Range("R5").Value = 1573.199526972405
$ws.Range("S5").Value = 0.001323095634797051
$ws.Range("T5").Value = 0.001323095634797051
$ws.Range("G6").Value = 47.86240033333333
$ws.Range("H6").Value = 143.587201
$ws.Range("I6").Value = 0.1228118231805696
$ws.Range("J6").Value = 0.1228118231805696
$ws.Range("M6").Value = 8.984181666666666
$ws.Range("N6").Value = 26.952545
$ws.Range("O6").Value = 0.02650225084781963
$ws.Range("P6").Value = 0.02650225084781963
$ws.Range("Q6").Value = 430.0044995973939
$ws.Range("R6").Value = 3870.040496376545
$ws.Range("S6").Value = 0.003254789745009526
$ws.Range("T6").Value = 0.003254789745009526
$ws.Range("I7").Value = 0.04786922362394307
$ws.Range("J7").Value = 0.04786922362394307
$ws.Range("M7").Value = 266.9240163333333
$ws.Range("N7").Value = 800.7720489999999
$ws.Range("O7").Value = 0.7873936103073201
$ws.Range("P7").Value = 0.78739361030732
$ws.Range("Q7").Value = 4979.644406346491
$ws.Range("R7").Value = 44816.7996571184
$ws.Range("S7").Value = 0.037691920811865
$ws.Range("T7").Value = 0.03769192081186499
$ws.Range("I8").Value = 0.04786922362394307
$ws.Range("J8").Value = 0.04786922362394307
$ws.Range("O8").Value = 0.1525285446808506
$ws.Range("P8").Value = 0.1525285446808506
$ws.Range("S8").Value = 0.007301423014362231
$ws.Range("T8").Value = 0.007301423014362231
$ws.Range("I9").Value = 0.04786922362394307
$ws.Range("J9").Value = 0.04786922362394307
$ws.Range("M9").Value = 7.729888333333332
$ws.Range("N9").Value = 23.189665
$ws.Range("O9").Value = 0.02280223700236483
$ws.Range("P9").Value = 0.02280223700236483
$ws.Range("Q9").Value = 144.2061891976689
$ws.Range("R9").Value = 1297.85570277902
$ws.Range("S9").Value = 0.001091525382192351
$ws.Range("T9").Value = 0.001091525382192351
$ws.Range("I10").Value = 0.04786922362394307
$ws.Range("J10").Value = 0.04786922362394307
$ws.Range("M10").Value = 3.652135
$ws.Range("N10").Value = 10.956405
$ws.Range("O10").Value = 0.01077335716164485
$ws.Range("P10").Value = 0.01077335716164485
$ws.Range("Q10").Value = 68.13299857312667
$ws.Range("R10").Value = 613.19698715814
$ws.Range("S10").Value = 0.0005157122431513862
$ws.Range("T10").Value = 0.0005157122431513862
$ws.Range("I11").Value = 0.04786922362394307
$ws.Range("J11").Value = 0.04786922362394307
$ws.Range("M11").Value = 8.984181666666666
$ws.Range("N11").Value = 26.952545
$ws.Range("O11").Value = 0.02650225084781963
$ws.Range("P11").Value = 0.02650225084781963
$ws.Range("Q11").Value = 167.6058625093845
$ws.Range("R11").Value = 1508.45276258446
$ws.Range("S11").Value = 0.001268642172372113
$ws.Range("T11").Value = 0.001268642172372113
$ws.Range("G12").Value = 171.0598806666667
$ws.Range("H12").Value = 513.1796420000001
$ws.Range("I12").Value = 0.4389285884413335
$ws.Range("J12").Value = 0.4389285884413335
$ws.Range("M12").Value = 266.9240163333333
$ws.Range("N12").Value = 800.7720489999999
$ws.Range("O12").Value = 0.7873936103073201
$ws.Range("P12").Value = 0.78739361030732
$ws.Range("Q12").Value = 45659.99038104739
$ws.Range("R12").Value = 410939.9134294264
$ws.Range("S12").Value = 0.3456095659199174
$ws.Range("T12").Value = 0.3456095659199174
$ws.Range("G13").Value = 171.0598806666667
$ws.Range("H13").Value = 513.1796420000001
$ws.Range("I13").Value = 0.4389285884413335
$ws.Range("J13").Value = 0.4389285884413335
$ws.Range("O13").Value = 0.1525285446808506
$ws.Range("P13").Value = 0.1525285446808506
$ws.Range("Q13").Value = 8844.943357166139
$ws.Range("R13").Value = 79604.49021449526
$ws.Range("S13").Value = 0.06694913881377662
$ws.Range("T13").Value = 0.06694913881377662
$ws.Range("G14").Value = 171.0598806666667
$ws.Range("H14").Value = 513.1796420000001
$ws.Range("I14").Value = 0.4389285884413335
$ws.Range("J14").Value = 0.4389285884413335
$ws.Range("M14").Value = 7.729888333333332
$ws.Range("N14").Value = 23.189665
$ws.Range("O14").Value = 0.02280223700236483
$ws.Range("P14").Value = 0.02280223700236483
$ws.Range("Q14").Value = 1322.273775866659
$ws.Range("R14").Value = 11900.46398279993
$ws.Range("S14").Value = 0.01000855370075274
$ws.Range("T14").Value = 0.01000855370075274
$ws.Range("G15").Value = 171.0598806666667
$ws.Range("H15").Value = 513.1796420000001
$ws.Range("I15").Value = 0.4389285884413335
$ws.Range("J15").Value = 0.4389285884413335
$ws.Range("M15").Value = 3.652135
$ws.Range("N15").Value = 10.956405
$ws.Range("O15").Value = 0.01077335716164485
$ws.Range("P15").Value = 0.01077335716164485
$ws.Range("Q15").Value = 624.7337772785568
$ws.Range("R15").Value = 5622.603995507011
$ws.Range("S15").Value = 0.004728734451735107
$ws.Range("T15").Value = 0.004728734451735107
$ws.Range("G16").Value = 171.0598806666667
$ws.Range("H16").Value = 513.1796420000001
$ws.Range("I16").Value = 0.4389285884413335
$ws.Range("J16").Value = 0.4389285884413335
$ws.Range("M16").Value = 8.984181666666666
$ws.Range("N16").Value = 26.952545
$ws.Range("O16").Value = 0.02650225084781963
$ws.Range("P16").Value = 0.02650225084781963
$ws.Range("Q16").Value = 1536.833043787655
$ws.Range("R16").Value = 13831.49739408889
$ws.Range("S16").Value = 0.0116325955551516
$ws.Range("T16").Value = 0.0116325955551516
$ws.Range("G17").Value = 12.628047
$ws.Range("H17").Value = 37.884141
$ws.Range("I17").Value = 0.0324027517316099
$ws.Range("J17").Value = 0.0324027517316099
$ws.Range("M17").Value = 266.9240163333333
$ws.Range("N17").Value = 800.7720489999999
$ws.Range("O17").Value = 0.7873936103073201
$ws.Range("P17").Value = 0.78739361030732
$ws.Range("Q17").Value = 3370.729023686101
$ws.Range("R17").Value = 30336.56121317491
$ws.Range("S17").Value = 0.02551371966984409
$ws.Range("T17").Value = 0.02551371966984408
$ws.Range("G18").Value = 12.628047
$ws.Range("H18").Value = 37.884141
$ws.Range("I18").Value = 0.0324027517316099
$ws.Range("J18").Value = 0.0324027517316099
$ws.Range("O18").Value = 0.1525285446808506
$ws.Range("P18").Value = 0.1525285446808506
$ws.Range("Q18").Value = 652.954743048625
$ws.Range("R18").Value = 5876.592687437625
$ws.Range("S18").Value = 0.00494234456527737
$ws.Range("T18").Value = 0.00494234456527737
$ws.Range("G19").Value = 12.628047
$ws.Range("H19").Value = 37.884141
$ws.Range("I19").Value = 0.0324027517316099
$ws.Range("J19").Value = 0.0324027517316099
$ws.Range("M19").Value = 7.729888333333332
$ws.Range("N19").Value = 23.189665
$ws.Range("O19").Value = 0.02280223700236483
$ws.Range("P19").Value = 0.02280223700236483
$ws.Range("Q19").Value = 97.61339317808499
$ws.Range("R19").Value = 878.5205386027649
$ws.Range("S19").Value = 0.0007388552245129561
$ws.Range("T19").Value = 0.0007388552245129561
$ws.Range("G20").Value = 12.628047
$ws.Range("H20").Value = 37.884141
$ws.Range("I20").Value = 0.0324027517316099
$ws.Range("J20").Value = 0.0324027517316099
$ws.Range("M20").Value = 3.652135
$ws.Range("N20").Value = 10.956405
$ws.Range("O20").Value = 0.01077335716164485
$ws.Range("P20").Value = 0.01077335716164485
$ws.Range("Q20").Value = 46.119332430345
$ws.Range("R20").Value = 415.073991873105
$ws.Range("S20").Value = 0.0003490864174247397
$ws.Range("T20").Value = 0.0003490864174247397
$ws.Range("G21").Value = 12.628047
$ws.Range("H21").Value = 37.884141
$ws.Range("I21").Value = 0.0324027517316099
$ws.Range("J21").Value = 0.0324027517316099
$ws.Range("M21").Value = 8.984181666666666
$ws.Range("N21").Value = 26.952545
$ws.Range("O21").Value = 0.02650225084781963
$ws.Range("P21").Value = 0.02650225084781963
$ws.Range("Q21").Value = 113.452668343205
$ws.Range("R21").Value = 1021.074015088845
$ws.Range("S21").Value = 0.0008587458545507474
$ws.Range("T21").Value = 0.0008587458545507474
$ws.Range("G22").Value = 139.5154473333333
$ws.Range("H22").Value = 418.546342
$ws.Range("I22").Value = 0.3579876130225438
$ws.Range("J22").Value = 0.3579876130225438
$ws.Range("M22").Value = 266.9240163333333
$ws.Range("N22").Value = 800.7720489999999
$ws.Range("O22").Value = 0.7873936103073201
$ws.Range("P22").Value = 0.78739361030732
$ws.Range("Q22").Value = 37240.02354275497
$ws.Range("R22").Value = 335160.2118847947
$ws.Range("S22").Value = 0.2818771590631206
$ws.Range("T22").Value = 0.2818771590631206
$ws.Range("G23").Value = 139.5154473333333
$ws.Range("H23").Value = 418.546342
$ws.Range("I23").Value = 0.3579876130225438
$ws.Range("J23").Value = 0.3579876130225438
$ws.Range("O23").Value = 0.1525285446808506
$ws.Range("P23").Value = 0.1525285446808506
$ws.Range("Q23").Value = 7213.884543259194
$ws.Range("R23").Value = 64924.96088933275
$ws.Range("S23").Value = 0.05460332962810013
$ws.Range("T23").Value = 0.05460332962810013
$ws.Range("G24").Value = 139.5154473333333
$ws.Range("H24").Value = 418.546342
$ws.Range("I24").Value = 0.3579876130225438
$ws.Range("J24").Value = 0.3579876130225438
$ws.Range("M24").Value = 7.729888333333332
$ws.Range("N24").Value = 23.189665
$ws.Range("O24").Value = 0.02280223700236483
$ws.Range("P24").Value = 0.02280223700236483
$ws.Range("Q24").Value = 1078.438828661714
$ws.Range("R24").Value = 9705.949457955428
$ws.Range("S24").Value = 0.008162918396050909
$ws.Range("T24").Value = 0.008162918396050909
$ws.Range("G25").Value = 139.5154473333333
$ws.Range("H25").Value = 418.546342
$ws.Range("I25").Value = 0.3579876130225438
$ws.Range("J25").Value = 0.3579876130225438
$ws.Range("M25").Value = 3.652135
$ws.Range("N25").Value = 10.956405
$ws.Range("O25").Value = 0.01077335716164485
$ws.Range("P25").Value = 0.01077335716164485
$ws.Range("Q25").Value = 509.5292482467233
$ws.Range("R25").Value = 4585.76323422051
$ws.Range("S25").Value = 0.003856728414536569
$ws.Range("T25").Value = 0.003856728414536569
$ws.Range("G26").Value = 139.5154473333333
$ws.Range("H26").Value = 418.546342
$ws.Range("I26").Value = 0.3579876130225438
$ws.Range("J26").Value = 0.3579876130225438
$ws.Range("M26").Value = 8.984181666666666
$ws.Range("N26").Value = 26.952545
$ws.Range("O26").Value = 0.02650225084781963
$ws.Range("P26").Value = 0.02650225084781963
$ws.Range("Q26").Value = 1253.432124148932
$ws.Range("R26").Value = 11280.88911734039
$ws.Range("S26").Value = 0.009487477520735638
$ws.Range("T26").Value = 0.009487477520735638

Write-Output "Updated 278 cells"